$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: 01/01/2012 -> 01/01/2021
# (Assigning a date-like literal directly to .Value would make Excel
#  auto-convert it to a date serial number and rewrite the cell's
#  number format. Instead, write it as a text formula, then convert that
#  formula to a plain value in place via Copy + PasteSpecial(values) so
#  the shared string / cell style stay exactly as they were.)
$dateText = "01/01/2021"

$b8 = $ws.Range("B8")
$b8.Formula = '="' + $dateText + '"'
$b8.Copy() | Out-Null
$b8.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues) | Out-Null

$c8 = $ws.Range("C8")
$c8.Formula = '="' + $dateText + '"'
$c8.Copy() | Out-Null
$c8.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues) | Out-Null

$excel.CutCopyMode = $false

# Docentes responsáveis:
$ws.Range("B13").Value = "8188658 - Maria Auxiliadora Motta Barreto"
$ws.Range("C13").Value = "8188658 - Maria Auxiliadora Motta Barreto"

# Programa resumido: (paragraphs collapsed into one line)
$programaResumido = "1. As necessidades das empresas modernas.2. Os novos modelos de produção. 3. As escolas de organização do trabalho. 4. Princípios sócio-técnicos de planejamento do trabalho. 5. Trabalho em grupo6. Organização por processos.7. Mudanças organizacionais."
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Programa: (paragraphs collapsed into one line + content changes)
$programa = "1. As necessidades das empresas modernas: Organização, produtividade, qualidade, flexibilidade e competitividade.2. Os novos modelos de produção: a ""revolução"" contemporânea nas fábricas.3. As escolas de organização do trabalho: administração científica, escola clássica (taylorismo/fordismo), escola de relações humanas (enriquecimento de cargos), escola sócio-técnica (grupos semi-autônomos), da contingência4. Princípios sócio-técnicos de planejamento do trabalho: metodologia original de projeto organizacional de Tavistock.5. Trabalho em grupo: tipos, casos. Grupos abertos e grupos fechados. Relações de fronteira (produção-manutenção, qualidade e planejamento, dentre outros).6. Organização por processos.7. Implantação de mudanças organizacionais"
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Método:
$ws.Range("B19").Value = "Apresentação de seminário e  prova escrita"
$ws.Range("C19").Value = "Apresentação de seminário e  prova escrita"

# Critério:
$ws.Range("B20").Value = "Média aritmética das duas  atividades avaliativas."
$ws.Range("C20").Value = "Média aritmética das duas  atividades avaliativas."

# Bibliografia: (paragraphs collapsed into one line)
$bibliografia = "FLEURY, Afonso C.C.; VARGAS, N. Organização do trabalho. São Paulo: Atlas, 1983.HELOANI, Roberto. Organização do Trabalho e Administração: Uma Visão Multidisciplinar. São Paulo: Cortez, 1994.MARX, Roberto. Trabalho em grupo e autonomia como instrumentos da competição. São Paulo: Atlas, 1998.SALERNO, Mario S. Projeto organizacional de produção integrada e flexível. São Paulo: EPUSP/PRO, 1998.SOUZA LIMA, J. C. Relações entre empresas da cadeia e questões contemporâneas de organização do trabalho. Itu: Ottoni, 2006.TAYLOR, Frederick W. Princípios de administração científica. São Paulo: Atlas, 1976.VARGAS, Nilton. Organização do trabalho e capital. Rio de Janeiro: COPPE/UFRJ/PEP, 1979.WOOMACK, James P.; JONES, Daniel T.; ROOS, Daniel. A máquina que mudou o mundo. Rio de Janeiro: Campus, 1992."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
